$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells are stored as TEXT in this sheet (e.g. "57.280.88", "0.562")
# Force text format before assigning so Excel does not auto-convert to a number,
# then reset the cell style back to Normal so no stray style index is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.280.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.553.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.23%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.564.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0997"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.94%  "

$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.009.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.307.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.571.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.398"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.673.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0740"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.60%  "

$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.42%  "

$ws.Range("E37").Value = "  -4.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.862"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.823"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "268.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0950"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.580"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.959.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.90%  "

# Row 34: EthereumClassic -> Monero
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.37%  "

# Row 35: Monero -> EthereumClassic
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.62%  "

# Row 51: InjectiveProtocol -> RenderToken
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.30%  "
